$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 377
$ws.Range("I2").Value = 941
$ws.Range("J2").Value = 3849
$ws.Range("K2").Value = 18
$ws.Range("L2").Value = 1083
$ws.Range("M2").Value = 55
$ws.Range("N2").Value = 663
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 11
$ws.Range("Q2").Value = 7
$ws.Range("R2").Value = 48
$ws.Range("S2").Value = 464
$ws.Range("T2").Value = 666
$ws.Range("U2").Value = 41
$ws.Range("V2").Value = 5864
$ws.Range("X2").Value = 5971
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 81
$ws.Range("AA2").Value = 33
